# feat: add 2022-Q1 data
#
# The existing "总计" sheet (summary table) becomes the new "2022-Q1" sheet
# (keeping its original sheet position/id) and is repopulated with that
# quarter's fund-holding data, matching the other quarterly sheets. A
# duplicate of the original "总计" sheet is inserted right after it and
# becomes the new "总计" summary, with a new leading row for 2022-Q1 on top
# of the previous summary rows.

$wb = $excel.ActiveWorkbook
$styleDonor = $wb.Worksheets.Item("2021-Q4")

# --- Step 0: duplicate 总计 (before touching it) so the copy inherits all
# sheet-level properties (sheetPr/pageMargins/...) verbatim, and lands right
# after the original in tab order. ------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Copy($null, $q1)
$totals = $wb.Worksheets.Item($q1.Index + 1)

# --- Step 1: turn the original 总计 sheet into the 2022-Q1 sheet -----------
$q1.Name = "2022-Q1"

# Clear out the old "日期 / 持有数量(只) / 持有市值(亿元)" summary table
# (values + formatting) before rebuilding it as a fund-holding table.
$q1.Cells.Clear()

# Header row.
$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"

# Data row. 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值 are stored
# as text (like every other quarter sheet), 仓位排名 is numeric.
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("A2").Value2 = 0
$q1.Range("B2").Value2 = "159792"
$q1.Range("C2").Value2 = "富国中证港股通互联网ETF"
$q1.Range("D2").Value2 = "2.76"
$q1.Range("E2").Value2 = "99.00"
$q1.Range("F2").Value2 = "2.60"
$q1.Range("G2").Value2 = "0.0718"
$q1.Range("H2").Value2 = 9

# Match the bold/bordered header + row-index-column formatting the other
# quarter sheets use (copy it verbatim instead of rebuilding it by hand).
$styleDonor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleDonor.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

# --- Step 2: rename the duplicate to 总计 and refresh its rows -------------
$totals.Name = "总计"

$totals.Range("A2").Value2 = 0
$totals.Range("B2").Value2 = "2022-Q1"
$totals.Range("C2").Value2 = 1
$totals.Range("D2").Value2 = 0.07000000000000001

$totals.Range("A3").Value2 = 1
$totals.Range("B3").Value2 = "2021-Q4"
$totals.Range("C3").Value2 = 6
$totals.Range("D3").Value2 = 1.6

$totals.Range("A4").Value2 = 2
$totals.Range("B4").Value2 = "2021-Q3"
$totals.Range("C4").Value2 = 10
$totals.Range("D4").Value2 = 4.19

$totals.Range("A5").Value2 = 3
$totals.Range("B5").Value2 = "2021-Q2"
$totals.Range("C5").Value2 = 10
$totals.Range("D5").Value2 = 4.22

$totals.Range("A6").Value2 = 4
$totals.Range("B6").Value2 = "2021-Q1"
$totals.Range("C6").Value2 = 7
$totals.Range("D6").Value2 = 2.95

# Row 6 is new (the old summary table only went down to row 5) so its A-cell
# needs the same row-index style as A2:A5 copied onto it explicitly.
$styleDonor.Range("A2").Copy()
$totals.Range("A6").PasteSpecial(-4122)
